# Regenerate merged AHB files
#
# 1) Rename the header columns from the "_old"/"_new" naming convention to the
#    concrete format-version names "_FV2304"/"_FV2310".
# 2) Turn the used range (A1:U58) into a real Excel Table (ListObject) with an
#    AutoFilter, using the renamed headers.
# 3) Freeze the header row (split after row 1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$suffixMap = @{
    "_old" = "_FV2304"
    "_new" = "_FV2310"
}

for ($i = 1; $i -le 21; $i++) {
    $cell = $ws.Cells.Item(1, $i)
    $text = [string]$cell.Value2
    if ($text -like "*_old") {
        $cell.Value = ($text -replace "_old$", "_FV2304")
    } elseif ($text -like "*_new") {
        $cell.Value = ($text -replace "_new$", "_FV2310")
    }
}

$tableRange = $ws.Range("A1:U58")
$table = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$table.Name = "Table1"

# Freeze the header row: put the selection on the first row below the
# header, then turn freeze panes on (matches Excel's "Freeze Top Row").
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
